# Validate fullname contains digit: append "0" to the Fullname value in C2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Hoang Viet Bach0"

# Move selection to D2 (matches recorded active cell after the edit)
$ws.Range("D2").Select()
